$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Docentes responsaveis" value row (old row 13, no label in column A)
# is being removed entirely; everything below shifts up one row.
$ws.Rows(13).Delete()

# After the shift, several cells further down were re-populated with
# different (shifted/duplicated) text values per the target revision.
$ws.Range("B10").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C10").Value = "11079086 - Herlandí de Souza Andrade"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2021"
$ws.Range("C15").Value = "01/01/2021"

$ws.Range("B18").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C18").Value = "11079086 - Herlandí de Souza Andrade"

$ws.Range("B19").Value = "Aulas expositivas teóricas, aulas práticas, aulas de exercícios."
$ws.Range("C19").Value = "Aulas expositivas teóricas, aulas práticas, aulas de exercícios."

$ws.Range("B20").Value = "Média Aritmética das atividades avaliativas realizadas."
$ws.Range("C20").Value = "Média Aritmética das atividades avaliativas realizadas."

$ws.Range("B21").Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."
$ws.Range("C21").Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."
